# Applies the "extracted data from 10.1016/j.jmst.2021.10.019" commit:
#  - adds new material/property rows (69-78) to Sheet1
#  - cell write order below is chosen to reproduce the upstream sharedStrings.xml
#    append order (new strings 146-163) as closely as possible
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Phase 1: new text values, in the exact order first typed (drives shared-string table order) ---
$ws.Range("B69").Value = "Fe40Mn14Ni10Cr10Al15C1"
$ws.Range("D69").Value = "BM"
$ws.Range("E69").Value = "ball milling for 40h until no elemental peaks were detected; 8nm crystalltie size with significant lattice strain"
$ws.Range("D70").Value = "BM+SPS"
$ws.Range("E70").Value = "ball milling for 40h until no elemental peaks were detected; SPS at 1273K under 50MPa"
$ws.Range("C70").Value = "FCC+B2+carbide+silicide"
$ws.Range("N69").Value = "10.1016/j.jallcom.2020.155013"
$ws.Range("B71").Value = "Fe2 Cr Ni Si0.3 Al0.28"
$ws.Range("D71").Value = "AAM"
$ws.Range("C71").Value = "BCC+FCC"
$ws.Range("E71").Value = "BCC majority"
$ws.Range("N71").Value = "10.1016/j.matlet.2023.134447"
$ws.Range("B74").Value = "Fe72.4 Co13.9 Cr10.4 Mn2.7 B0.34"
$ws.Range("E74").Value = "magnetron sputtering deposition; micropilars were tested"
$ws.Range("F76").Value = "UCS"
$ws.Range("F75").Value = "minimum compressive ductility"
$ws.Range("E77").Value = "magnetron sputtering deposition"
$ws.Range("N74").Value = "10.1016/j.jmst.2021.10.019"

# --- Phase 2: remaining text cells that reuse already-existing shared strings ---
$ws.Range("C69").Value = "BCC"
$ws.Range("B70").Value = "Fe40Mn14Ni10Cr10Al15C1"
$ws.Range("F70").Value = "hardness"
$ws.Range("G70").Value = "EXP"
$ws.Range("L70").Value = "Pa"
$ws.Range("N70").Value = "10.1016/j.jallcom.2020.155013"
$ws.Range("F71").Value = "tensile yield stress"
$ws.Range("G71").Value = "EXP"
$ws.Range("L71").Value = "Pa"
$ws.Range("B72").Value = "Fe2 Cr Ni Si0.3 Al0.28"
$ws.Range("C72").Value = "BCC+FCC"
$ws.Range("D72").Value = "AAM"
$ws.Range("E72").Value = "BCC majority"
$ws.Range("F72").Value = "UTS"
$ws.Range("G72").Value = "EXP"
$ws.Range("L72").Value = "Pa"
$ws.Range("N72").Value = "10.1016/j.matlet.2023.134447"
$ws.Range("B73").Value = "Fe2 Cr Ni Si0.3 Al0.28"
$ws.Range("C73").Value = "BCC+FCC"
$ws.Range("D73").Value = "AAM"
$ws.Range("E73").Value = "BCC majority"
$ws.Range("F73").Value = "tensile ductility"
$ws.Range("G73").Value = "EXP"
$ws.Range("L73").Value = "%"
$ws.Range("N73").Value = "10.1016/j.matlet.2023.134447"
$ws.Range("C74").Value = "BCC"
$ws.Range("D74").Value = "SD"
$ws.Range("F74").Value = "compressive yield strength"
$ws.Range("G74").Value = "EXP"
$ws.Range("L74").Value = "Pa"
$ws.Range("B75").Value = "Fe72.4 Co13.9 Cr10.4 Mn2.7 B0.34"
$ws.Range("C75").Value = "BCC"
$ws.Range("D75").Value = "SD"
$ws.Range("E75").Value = "magnetron sputtering deposition; micropilars were tested"
$ws.Range("G75").Value = "EXP"
$ws.Range("L75").Value = "%"
$ws.Range("N75").Value = "10.1016/j.jmst.2021.10.019"
$ws.Range("B76").Value = "Fe72.4 Co13.9 Cr10.4 Mn2.7 B0.34"
$ws.Range("C76").Value = "BCC"
$ws.Range("D76").Value = "SD"
$ws.Range("E76").Value = "magnetron sputtering deposition; micropilars were tested"
$ws.Range("G76").Value = "EXP"
$ws.Range("L76").Value = "Pa"
$ws.Range("N76").Value = "10.1016/j.jmst.2021.10.019"
$ws.Range("B77").Value = "Fe72.4 Co13.9 Cr10.4 Mn2.7 B0.34"
$ws.Range("C77").Value = "BCC"
$ws.Range("D77").Value = "SD"
$ws.Range("F77").Value = "nanohardness"
$ws.Range("G77").Value = "EXP"
$ws.Range("L77").Value = "Pa"
$ws.Range("N77").Value = "10.1016/j.jmst.2021.10.019"
$ws.Range("B78").Value = "Fe72.4 Co13.9 Cr10.4 Mn2.7 B0.34"
$ws.Range("C78").Value = "BCC"
$ws.Range("D78").Value = "SD"
$ws.Range("E78").Value = "magnetron sputtering deposition"
$ws.Range("F78").Value = "reduced elastic modulus"
$ws.Range("G78").Value = "EXP"
$ws.Range("L78").Value = "Pa"
$ws.Range("N78").Value = "10.1016/j.jmst.2021.10.019"

# --- Phase 3: numeric + formula cells ---
$ws.Range("I70").Value = 298
$ws.Range("J70").Formula = "=P70*9807000"
$ws.Range("P70").Value = 596
$ws.Range("I71").Value = 298
$ws.Range("J71").Value = 633000000
$ws.Range("K71").Value = 49000000
$ws.Range("I72").Value = 298
$ws.Range("J72").Value = 1259000000
$ws.Range("K72").Value = 42000000
$ws.Range("I73").Value = 298
$ws.Range("J73").Value = 9.32
$ws.Range("K73").Value = 0.18
$ws.Range("I74").Value = 298
$ws.Range("J74").Value = 2920000000
$ws.Range("K74").Value = 0.36
$ws.Range("I75").Value = 298
$ws.Range("J75").Value = 13.7
$ws.Range("K75").Value = 1.9
$ws.Range("I76").Value = 298
$ws.Range("J76").Value = 3370000000
$ws.Range("K76").Value = 0.36
$ws.Range("I77").Value = 298
$ws.Range("J77").Value = 9300000000
$ws.Range("K77").Value = 400000000
$ws.Range("I78").Value = 298
$ws.Range("J78").Value = 293000000000
$ws.Range("K78").Value = 12000000000

# --- Update view state: scrolled down a bit further, selection moved to M80 ---
$ws.Activate()
$ws.Range("M80").Select() | Out-Null
